# Update generated output values (columns F/G) across all sheets
# per "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 198
$ws.Range("F6").Value = 795
$ws.Range("F7").Value = 90
$ws.Range("F8").Value = 10141
$ws.Range("F10").Value = 3498
$ws.Range("F12").Value = 2435
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 2787
$ws.Range("F17").Value = 2154
$ws.Range("F18").Value = 46
$ws.Range("F20").Value = 19
$ws.Range("F21").Value = 383
$ws.Range("F23").Value = 135
$ws.Range("F26").Value = 214
$ws.Range("F28").Value = 1311
$ws.Range("F30").Value = 1249
$ws.Range("F31").Value = 102
$ws.Range("F34").Value = 3180
$ws.Range("F35").Value = 3022
$ws.Range("F36").Value = 25
$ws.Range("F38").Value = 1031
$ws.Range("F39").Value = 393
$ws.Range("F41").Value = 1291
$ws.Range("F42").Value = 89
$ws.Range("F43").Value = 106
$ws.Range("F44").Value = 71

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("G15").Value = 266
$ws.Range("F16").Value = 175

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 741
$ws.Range("F5").Value = 1990

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F8").Value = 198
$ws.Range("F9").Value = 90
$ws.Range("F10").Value = 10141
$ws.Range("F12").Value = 3498
$ws.Range("F14").Value = 27
$ws.Range("F17").Value = 2154
$ws.Range("F18").Value = 46
$ws.Range("F20").Value = 19
$ws.Range("F21").Value = 135
$ws.Range("F24").Value = 214
$ws.Range("F26").Value = 1311
$ws.Range("F28").Value = 1249
$ws.Range("F33").Value = 3181
$ws.Range("F34").Value = 3022
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 1031
$ws.Range("F39").Value = 393
$ws.Range("G43").Value = 266
$ws.Range("F44").Value = 89
$ws.Range("F45").Value = 71
$ws.Range("F49").Value = 175
